# Refactor synthetic array /3
# Swap the "noir" (black) status entry for a "bleu" (blue) one:
#   - emoji square ⬛ -> book 📘
#   - emoji square 🟥 -> book 📕
#   - emoji square 🟧 -> book 📙
#   - emoji square 🟩 -> book 📗
#   - label "noir" -> "bleu"  (rouge / orange / vert untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "⬛" = "📘"
    "🟥" = "📕"
    "🟧" = "📙"
    "🟩" = "📗"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cellA = $ws.Range("A$r")
    $valA = $cellA.Value2
    if ($valA -ne $null -and $map.ContainsKey($valA)) {
        $cellA.Value2 = $map[$valA]
    }

    $cellB = $ws.Range("B$r")
    $valB = $cellB.Value2
    if ($valB -eq "noir") {
        $cellB.Value2 = "bleu"
    }
}
